$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for account 005651698 (JOAO) and 004452597 (LARA)
$ws.Range("A7:A8").EntireRow.Delete()
